# Updates cryptos list cell values to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.541.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.817.55'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.68%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.88'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.559'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.59'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +7.15%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0694'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0951'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.079.90'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.33'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.809.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.644'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.581.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.33'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.06'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '247.23'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0801'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.33%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.21'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '172.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.63%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.70'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0530'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.25'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.84'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.61'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.418.58'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.673'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.35%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '85.81'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.953'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.56%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.84'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0527'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.57%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.980.61'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.68'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.04%  '

Write-Output "Applied 91 cell updates"
